$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix DIVA_ID value in row 2 (was 1, should be 466)
$ws.Range("A2").Value = 466

# Duplicate row 2 into rows 3 and 4 (copy values + formatting)
$ws.Range("A2:AU2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial() | Out-Null
$ws.Range("A2:AU2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial() | Out-Null

# Update the active selection to match the post-edit cursor position
$ws.Range("A4").Select() | Out-Null
